$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so that numeric-looking
# values (e.g. "307.85") are written as strings, matching the source data
# which stores all Price values as inline/shared strings rather than numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.248.29'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '1.902.17'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '307.85'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.5216'
$ws.Range("E7").Value = '  +0.63%  '
$ws.Range("D8").Value = '0.3776'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '0.07289'
$ws.Range("E9").Value = '  +0.98%  '
$ws.Range("D10").Value = '21.23'
$ws.Range("E10").Value = '  +0.16%  '
$ws.Range("D11").Value = '0.9031'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '0.08294'
$ws.Range("E12").Value = '  +8.43%  '
$ws.Range("D13").Value = '96.64'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("D14").Value = '1.893.45'
$ws.Range("E14").Value = '  +2.75%  '
$ws.Range("D15").Value = '5.291'
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '0.000008630'
$ws.Range("E17").Value = '  +1.38%  '
$ws.Range("D18").Value = '14.57'
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("D19").Value = '0.9995'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = '27.278.26'
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").Value = '5.094'
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("D22").Value = '2.141.71'
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("D23").Value = '10.68'
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("D24").Value = '6.438'
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("D25").Value = '2.313'
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("D26").Value = '147.31'
$ws.Range("E26").Value = '  +1.20%  '
$ws.Range("D27").Value = '1.752'
$ws.Range("E27").Value = '  +1.62%  '
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").Value = '115.36'
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("D30").Value = '4.841'
$ws.Range("E30").Value = '  +0.84%  '
$ws.Range("D31").Value = '4.912'
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("D32").Value = '0.09252'
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").Value = '0.05068'
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("D34").Value = '0.7999'
$ws.Range("E34").Value = '  +3.64%  '
$ws.Range("D35").Value = '1.238'
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("D36").Value = '3.438'
$ws.Range("E36").Value = '  +4.73%  '
$ws.Range("D37").Value = '2.955'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").Value = '2.590'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = '0.5723'
$ws.Range("E39").Value = '  +1.92%  '
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").Value = '1.082'
$ws.Range("E41").Value = '  +1.03%  '
$ws.Range("D42").Value = '9.021'
$ws.Range("E42").Value = '  -0.42%  '
$ws.Range("D43").Value = '6.581'
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("D44").Value = '116.08'
$ws.Range("E44").Value = '  -2.13%  '
$ws.Range("D45").Value = '0.1519'
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").Value = '0.4876'
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").Value = '0.9997'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = '10.09'
$ws.Range("E48").Value = '  -0.70%  '
$ws.Range("D49").Value = '1.628'
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("D50").Value = '38.03'
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("D51").Value = '63.86'
$ws.Range("E51").Value = '  +0.01%  '

# Restore the default "Normal" style on the Price column so no stray
# number-format style is left attached to these cells.
$priceRange.Style = "Normal"
